# Generate Report for Handoff
# Updates status text "In Translation" -> "Ready for handoff" and
# refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps, also widening the status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text updates --------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamp updates -----------------------------------------------------
$wsOverview.Range("G2").Value = "2016-08-30 13:03:23"
$wsDeDe.Range("H2").Value = "2016-08-30 13:03:23"
$wsZhCn.Range("H2").Value = "2016-08-30 13:03:19"

# --- Column width updates ---------------------------------------------------
# (target character width ~17.22; engine snaps ColumnWidth to the nearest
# pixel grid, so 16.3333 is the input that lands closest to that target)
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
